$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '67.464.82'
$ws.Range("E2").Value = '  +2.47%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.587.73'
$ws.Range("E3").Value = '  +1.35%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = '''196.44'
$ws.Range("E5").Value = '  +7.06%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '''555.78'
$ws.Range("E6").Value = '  -4.48%  '

$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.582.64'
$ws.Range("E7").Value = '  +1.46%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '''0.609'
$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  -0.30%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '''0.665'
$ws.Range("E10").Value = '  +0.86%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.150'
$ws.Range("E11").Value = '  +5.27%  '

$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").Value = '''55.90'
$ws.Range("E12").Value = '  +5.15%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '''0.0000286'
$ws.Range("E13").Value = '  +14.82%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''9.86'
$ws.Range("E14").Value = '  +2.19%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.165.01'
$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.598.57'
$ws.Range("E16").Value = '  +1.61%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '''0.126'
$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '67.368.56'
$ws.Range("E18").Value = '  +2.69%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '''18.40'
$ws.Range("E19").Value = '  +1.63%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''12.15'
$ws.Range("E20").Value = '  +1.25%  '

$ws.Range("B21").Value = 'Polygon'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D21").Value = '''1.07'
$ws.Range("E21").Value = '  +2.84%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''393.17'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("B23").Value = 'RenderToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D23").Value = '''13.03'
$ws.Range("E23").Value = '  +26.30%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '''4.11'
$ws.Range("E24").Value = '  -3.93%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''84.12'
$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '''2.91'
$ws.Range("E26").Value = '  +2.11%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '''12.31'
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").Value = '''6.09'
$ws.Range("E28").Value = '  +1.20%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''3.79'
$ws.Range("E29").Value = '  +8.16%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '''8.20'
$ws.Range("E30").Value = '  +21.71%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '''8.94'
$ws.Range("E31").Value = '  +1.32%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '''31.26'
$ws.Range("E32").Value = '  +2.05%  '

$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '''676.17'
$ws.Range("E33").Value = '  +11.22%  '

$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").Value = '''12.06'
$ws.Range("E34").Value = '  +1.02%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.114'
$ws.Range("E35").Value = '  +2.78%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '''63.47'
$ws.Range("E36").Value = '  +2.59%  '

$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").Value = '''42.14'
$ws.Range("E37").Value = '  +4.06%  '

$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").Value = '''0.431'
$ws.Range("E38").Value = '  +17.32%  '

$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0762'
$ws.Range("E40").Value = '  +4.11%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.137'
$ws.Range("E41").Value = '  +6.62%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.203.36'
$ws.Range("E42").Value = '  +11.31%  '

$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Value = '''2.83'
$ws.Range("E43").Value = '  +17.83%  '

$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").Value = '''3.12'
$ws.Range("E44").Value = '  +14.07%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''0.999'
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '''3.00'
$ws.Range("E46").Value = '  +34.81%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0412'
$ws.Range("E47").Value = '  +1.79%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '''2.68'
$ws.Range("E48").Value = '  +10.02%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '''3.13'
$ws.Range("E49").Value = '  +2.95%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '''0.130'
$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '''8.56'
$ws.Range("E51").Value = '  +2.61%  '
